# Apply "Arthen's implementation for reduction in dividends" adjustments
# to the percentile income indicators table.
# Updates columns D (P99), F (Dif_P99_P50), H (Razao_P99_P50),
# I (Media_Top_0.1) and J (Razao_Top001_P50) for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Regime Atual
$ws.Range("D2").Value = 34958.9187444724
$ws.Range("F2").Value = 33230.29547920614
$ws.Range("H2").Value = 20.22356140109434
$ws.Range("I2").Value = 488104.8067994241
$ws.Range("J2").Value = 282.3662139733155

# Row 3 - Nova Proposta
$ws.Range("D3").Value = 34958.9187444724
$ws.Range("F3").Value = 33230.29547920614
$ws.Range("H3").Value = 20.22356140109434
$ws.Range("I3").Value = 465205.9881384756
$ws.Range("J3").Value = 269.1193607571976

# Row 4 - Nova c/ Aliq. Máxima
$ws.Range("D4").Value = 34853.059284889
$ws.Range("F4").Value = 33124.43601962274
$ws.Range("H4").Value = 20.16232222786877
$ws.Range("I4").Value = 458592.727236731
$ws.Range("J4").Value = 265.2936220698696
